$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric/Volume table updates (rows 2-50) ---
$ws.Range("D2").Value = "'58.903.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "'2.654.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'512.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'143.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "'0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("D9").Value = "'2.656.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("D10").Value = "'6.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("D11").Value = "'0.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").Value = "'0.336"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D14").Value = "'3.112.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "'58.892.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "'21.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "'2.652.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").Value = "'4.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").Value = "'343.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "'10.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "'6.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D24").Value = "'60.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("D26").Value = "'2.757.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("D27").Value = "'0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("D29").Value = "'0.0₃0803"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("D30").Value = "'7.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.26%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "'6.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.37%  "
$ws.Range("D35").Value = "'149.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.26%  "
$ws.Range("D37").Value = "'4.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("D39").Value = "'0.851"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("D40").Value = "'36.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").Value = "'3.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").Value = "'1.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'280.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").Value = "'0.614"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").Value = "'0.998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "'0.0983"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").Value = "'19.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").Value = "'0.0531"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").Value = "'10.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").Value = "'0.0228"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.65%  "

# --- Row 33/34 swap: PancakeSwap <-> EthereumClassic ---
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'18.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "

$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "'1.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.31%  "

# --- Row 51: Maker -> RenderToken ---
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'4.68"
$ws.Range("D51").Style = "Normal"
